$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Optimizer"); this shifts
# Optimizer..val_accuracy (old D..N) right by one into E..O.
$ws.Columns.Item(4).Insert()

# Header + data for the new "Description" column
$ws.Cells.Item(1, 4).Value = "Description"
$ws.Cells.Item(2, 4).Value = "pretrained  model weights, fc layer alone trained"
$ws.Cells.Item(3, 4).Value = "pretrained  model weights, fc layer alone trained"

# Row heights grew slightly to accommodate the new/rewrapped content
$ws.Rows.Item(1).RowHeight = 23.85
$ws.Rows.Item(2).RowHeight = 68.65
$ws.Rows.Item(3).RowHeight = 67.95

# Column widths were slightly re-balanced across the sheet
$ws.Columns.Item(1).ColumnWidth = 3.35034013605443
$ws.Columns.Item(2).ColumnWidth = 14.4217687074830
$ws.Columns.Item(3).ColumnWidth = 15.2329931972789
$ws.Columns.Item(4).ColumnWidth = 10.5595238095238
$ws.Columns.Item(5).ColumnWidth = 8.34523809523810
$ws.Columns.Item(6).ColumnWidth = 7.67176870748300
$ws.Columns.Item(7).ColumnWidth = 7.94217687074830
$ws.Columns.Item(8).ColumnWidth = 9.82993197278917
$ws.Columns.Item(9).ColumnWidth = 6.31972789115647
$ws.Columns.Item(10).ColumnWidth = 7.53401360544218
$ws.Columns.Item(11).ColumnWidth = 11.5850340136055
$ws.Columns.Item(12).ColumnWidth = 10.5085034013606
$ws.Columns.Item(13).ColumnWidth = 11.5850340136055
$ws.Columns.Item(14).ColumnWidth = 9.01870748299320
$ws.Columns.Item(15).ColumnWidth = 12.6666666666667

# Active selection moved to F3
$null = $ws.Range("F3").Select()
